$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Target journal proposals")

# Add the new suggested journal in row 12
$ws.Range("A12").Value = "Open Research Europe"
$ws.Range("F12").Value = "Andreas"

# L12 holds the "Info" hyperlink for the journal
$ws.Hyperlinks.Add($ws.Range("L12"), "https://open-research-europe.ec.europa.eu/") | Out-Null
$ws.Range("L12").Style = $ws.Range("L2").Style

$ws.Range("G12").Value = "Compliance with EU's Open Science policy"

# Update the active selection to match the final state recorded in the workbook
$ws.Range("G11").Select()
